$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("boosters")
$ws.Activate()

# Shift the existing "nicotine"/"volume" data from columns C:D to E:F so
# two new columns ("pg"/"vg") can be introduced at C:D. The column width
# metadata already defined for columns C/D is left untouched (it will now
# describe the new pg/vg columns); new width metadata is created for E/F.

# --- Header row ---
$ws.Range("E1").Value = "nicotine"
$ws.Range("F1").Value = "volume"
$ws.Range("E1:F1").Style = "Nadpis 2"

$ws.Range("C1").Value = "pg"
$ws.Range("D1").Value = "vg"
$ws.Range("C1:D1").Style = "Nadpis 2"

# --- Row 2 : Imperia / Dripper VPG 70/30 ---
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 10
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 70

# --- Row 3 : Imperia / Dripper VPG 70/30 ---
$ws.Range("E3").Value = 18
$ws.Range("F3").Value = 10
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 70

# --- Row 4 : Imperia / Dripper VPG 70/30 ---
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 10
$ws.Range("C4").Value = 30
$ws.Range("D4").Value = 70

# --- Row 5 : Imperia / Fifty VPG 50/50 ---
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = 10
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 50

# --- Row 6 : Imperia / Fifty VPG 50/50 ---
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 10
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 50

# --- Row 7 : Imperia / Fifty VPG 50/50 ---
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 10
$ws.Range("C7").Value = 50
$ws.Range("D7").Value = 50

# --- Row 8 : Imperia / Nico Base Fifty ---
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 10
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = 50

# Match column widths from the diff (new metadata for the shifted columns)
$ws.Range("E1").EntireColumn.ColumnWidth = 13.42578125
$ws.Range("F1").EntireColumn.ColumnWidth = 12.7109375

# Selection / active cell as seen in the final file
$ws.Range("D5").Select() | Out-Null
